$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# habitat_typenumber (F) and mfd_hab3 pollution status (P) per row
$habitatNumber = @{
    2 = "3110"; 3 = "3130"; 4 = "3130"; 5 = "3130"; 6 = "3130";
    7 = "3110"; 8 = "3110"; 9 = "3140"; 10 = "3110"; 11 = "3110";
    12 = "3110"; 13 = "3130"; 14 = "3110"; 15 = "3110"; 16 = "3110";
    17 = "3130"; 18 = "3130"; 19 = "3130"; 20 = "3130"; 21 = "3130";
    22 = "3110"; 23 = "3110"; 24 = "3110"; 25 = "3130"
}

$pollution = @{
    2 = "No pollution"; 3 = "Polluted"; 4 = "Polluted"; 5 = "Polluted"; 6 = "Polluted";
    7 = "No pollution"; 8 = "No pollution"; 9 = "Pollution results unknown"; 10 = "No pollution"; 11 = "No pollution";
    12 = "No pollution"; 13 = "Polluted"; 14 = "No pollution"; 15 = "No pollution"; 16 = "No pollution";
    17 = "Polluted"; 18 = "Polluted"; 19 = "Polluted"; 20 = "Polluted"; 21 = "Polluted";
    22 = "No pollution"; 23 = "No pollution"; 24 = "No pollution"; 25 = "Polluted"
}

# habitat_typenumber values look numeric ("3110", "3130", "3140") but must be
# stored as text, matching the source data. Format as text first so Excel
# doesn't coerce the values into numbers.
$ws.Range("F2:F25").NumberFormat = "@"

for ($row = 2; $row -le 25; $row++) {
    $ws.Range("F$row").Value = $habitatNumber[$row]
    $ws.Range("L$row").Value = "Subterranean"
    $ws.Range("N$row").Value = "Freshwater"
    $ws.Range("O$row").Value = "Groundwater"
    $ws.Range("P$row").Value = $pollution[$row]
}
